$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, shifting existing rows 37-140 down to 38-141
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly data entry
$ws.Cells.Item(37, 1).Value = 11
$ws.Cells.Item(37, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value = "Bíobío"
$ws.Cells.Item(37, 4).Value = 44980
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100101
$ws.Cells.Item(37, 8).Value = "Berries"
$ws.Cells.Item(37, 9).Value = 100101001
$ws.Cells.Item(37, 10).Value = "Arándano (blue)"
$ws.Cells.Item(37, 11).Value = "Sin especificar"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 150
$ws.Cells.Item(37, 14).Value = 3000
$ws.Cells.Item(37, 15).Value = 3500
$ws.Cells.Item(37, 16).Value = 3233
$ws.Cells.Item(37, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(37, 18).Value = "Región de Ñuble"
$ws.Cells.Item(37, 19).Value = 1616
$ws.Cells.Item(37, 20).Value = 2
